$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.479.73'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '2.214.72'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'" + '240.26'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = "'" + '0.610'
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("D7").Value = "'" + '74.81'
$ws.Range("E7").Value = '  +3.07%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = "'" + '0.599'
$ws.Range("E9").Value = '  +1.91%  '
$ws.Range("D10").Value = "'" + '41.20'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("D12").Value = "'" + '54.83'
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("E14").Value = '  -1.90%  '
$ws.Range("D15").Value = '2.549.00'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = "'" + '14.62'
$ws.Range("E16").Value = '  +3.21%  '
$ws.Range("D17").Value = '2.213.47'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = "'" + '0.799'
$ws.Range("E18").Value = '  -3.41%  '
$ws.Range("D19").Value = '42.370.34'
$ws.Range("E19").Value = '  +1.73%  '
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").Value = "'" + '70.64'
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("E22").Value = '  -3.53%  '
$ws.Range("D23").Value = "'" + '9.91'
$ws.Range("E23").Value = '  -10.05%  '
$ws.Range("D24").Value = "'" + '229.03'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").Value = "'" + '2.14'
$ws.Range("E25").Value = '  +5.93%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = "'" + '10.90'
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("D28").Value = "'" + '3.39'
$ws.Range("E28").Value = '  -6.45%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = "'" + '172.46'
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = "'" + '36.60'
$ws.Range("E31").Value = '  +19.66%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = "'" + '2.09'
$ws.Range("E32").Value = '  -4.94%  '
$ws.Range("D33").Value = "'" + '20.22'
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").Value = "'" + '0.0790'
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = "'" + '5.38'
$ws.Range("E35").Value = '  -1.80%  '
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").Value = "'" + '4.41'
$ws.Range("E38").Value = '  +4.35%  '
$ws.Range("D39").Value = "'" + '0.0320'
$ws.Range("E39").Value = '  +6.61%  '
$ws.Range("D40").Value = "'" + '12.44'
$ws.Range("E40").Value = '  -4.01%  '
$ws.Range("D41").Value = "'" + '2.13'
$ws.Range("D42").Value = "'" + '5.48'
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").Value = "'" + '60.39'
$ws.Range("E43").Value = '  -5.38%  '
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("D45").Value = "'" + '8.56'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").Value = "'" + '0.0990'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("D47").Value = "'" + '99.21'
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").Value = "'" + '2.27'
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  -2.60%  '
$ws.Range("D51").Value = "'" + '0.422'
$ws.Range("E51").Value = '  +15.66%  '
